# Precizejumi no Everitas par kodiem.
# Applies updates to Sheet1 per the author's edits: adds/updates indexing-notes
# cells (column C/E/F) and tweaks the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('E5').Value = 'Nē, neindeksē (Everita)'
$ws.Range('E11').Value = 'Nē, neindeksē (Everita). '
$ws.Range('F11').Value = 'Everita: Teorētiski tas arī ir teksts latviski, tāpēc, ja ir vēlme, to varētu indeksēt, bet adresē tad norādītu tikai, piem., JT1685 2P_g - bet tas tad būtu jāatrunā. Šobrīd mums jau ir 0. pants, ko mēs indeksējam, kas arī nav "oficiālais" Bībeles teksts. Bet šis @g{} ir nesvarīgs, tāpēc mēs varam turpināt to ignorēt. Piemēram @p{} ir daudz svarīgāks, ko vajadzētu indeksēt. sk. zemāk'
$ws.Range('C14').Value = 'manuprāt, brīvs apzīmējums :)'
$ws.Range('E17').Value = 'Nē (Everita)'
$ws.Range('C18').Value = 'Tas bija vecais, tikai vienā avotā izmantototais apzīmējums nodaļas nosaukumam augšā, kas atkārtojas katrā lapā. Šobrīd mēs to vairs neizmantojam.'
$ws.Range('E19').Value = 'Nē (Normunds). '
$ws.Range('F19').Value = 'Everita:  ja ir iespēja, tad mēs ļoti gribētu šīs Bībeles (!) piezīmes tekstu arī indeksēt, jo tas ir latviski. Adrese šādam vārdlietojumam būtu tāds pat kā Bībelei: nodaļa: pants+p , piem., kaut kāds (izgudrots piemērs) 1Sam 1:4p. Savukārt tur, kur piezīme ir ne Bībeles tekstā, to var ignorēt, jo tas ir tikai atsevišķš vārds vai skaņa, kurai turklāt grūti "pielikt" adresi (precīzu rindiņu), labākajā gadījumā tā ir tikai lappuses p: Manc1631_LVM 12p (izdomāts piemērs) -- bet tas lec ārā no pārējās šī avota adreses. Ne Bībeles tekstā šis ir mazsvarīgs.'
$ws.Range('C20').Value = 'Brīvs :-)'
$ws.Range('C24').Value = 'Brīvs :-)'
$ws.Range('C27').Value = 'Kaut kādi tukšumi? - Tas varētu būt izmantots xxx kā kaut kas nesalasāms, bet vai tas patiešām ir iekš {} ?'
$ws.Range('C28').Value = 'Brīvs :-)'
$ws.Range('E29').Value = 'Nē'

# Restore the author's final selection on the sheet.
$ws.Range("E19").Select()
